$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 322, shifting existing rows 322-351 down to 323-352
$ws.Rows.Item(322).Insert()

# Populate the newly inserted row 322 with data
$ws.Range("A322").Value = 9
$ws.Range("B322").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C322").Value = "Metropolitana"
$ws.Range("D322").Value = 44578
$ws.Range("E322").Value = 13
$ws.Range("F322").Value = 100112012
$ws.Range("G322").Value = "Espinaca"
$ws.Range("H322").Value = "Sin especificar"
$ws.Range("I322").Value = "Primera"
$ws.Range("J322").Value = 61
$ws.Range("K322").Value = 14000
$ws.Range("L322").Value = 15000
$ws.Range("M322").Value = 14508
$ws.Range("N322").Value = '$/cuna 10 kilos'
$ws.Range("O322").Value = "Provincia de Chacabuco"
$ws.Range("P322").Value = 1451
$ws.Range("Q322").Value = 10
$ws.Range("R322").Value = "Hortaliza"
